$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (sub_leaf_1): vein_angle_mean/median/std (C:E) and vein_cross_angle_1..4 (R:U)
$ws.Range("C2").Value = -6.782253339837299
$ws.Range("D2").Value = -0.20456564967241775
$ws.Range("E2").Value = 53.10851737952483
$ws.Range("R2").Value = -75.27386634281146
$ws.Range("S2").Value = -41.69146886217747
$ws.Range("T2").Value = 41.28233756283264
$ws.Range("U2").Value = 48.5539842828071

# Row 3 (sub_leaf_2)
$ws.Range("C3").Value = -0.30890568620245773
$ws.Range("D3").Value = 1.34241802464701
$ws.Range("E3").Value = 42.580041408615855
$ws.Range("R3").Value = -46.288270262894805
$ws.Range("S3").Value = 42.03259995738133
$ws.Range("T3").Value = 42.367811468790954
$ws.Range("U3").Value = -39.34776390808731

# Row 4 (sub_leaf_3)
$ws.Range("C4").Value = 33.92211654331513
$ws.Range("D4").Value = 51.48826472585644
$ws.Range("E4").Value = 44.875099521691986
$ws.Range("R4").Value = 74.2387529794948
$ws.Range("S4").Value = -41.52681625794717
$ws.Range("T4").Value = 59.233794727365634
$ws.Range("U4").Value = 43.74273472434724

# Row 5 (sub_leaf_4)
$ws.Range("C5").Value = 21.996248083321447
$ws.Range("D5").Value = 46.805060137145375
$ws.Range("E5").Value = 48.1850912717703
$ws.Range("R5").Value = 55.59021296997648
$ws.Range("S5").Value = -61.215340910981425
$ws.Range("T5").Value = 45.47563306716019
$ws.Range("U5").Value = 48.134487207130555

# Row 6 (sub_leaf_5)
$ws.Range("C6").Value = -7.193035369557796
$ws.Range("D6").Value = -4.964945230582462
$ws.Range("E6").Value = 46.50053412110406
$ws.Range("R6").Value = -63.478621269151944
$ws.Range("S6").Value = 44.63637025208568
$ws.Range("T6").Value = -42.2805951358408
$ws.Range("U6").Value = 32.350704674675875

# Row 7 (sub_leaf_6) - only 3 cross angles (R,S,T)
$ws.Range("C7").Value = -35.657918805343776
$ws.Range("D7").Value = -64.12503882579094
$ws.Range("E7").Value = 50.77754265668545
$ws.Range("S7").Value = 35.669910319068116
$ws.Range("T7").Value = -64.12503882579094

# Row 8 (sub_leaf_7) - 5 cross angles (R,S,T,U,V)
$ws.Range("C8").Value = 6.1227293717111895
$ws.Range("D8").Value = 39.48392088041287
$ws.Range("E8").Value = 43.986563962723174
$ws.Range("R8").Value = -49.986510396819845
$ws.Range("S8").Value = -45.37203067897739
$ws.Range("T8").Value = 39.48392088041287
$ws.Range("U8").Value = 44.85123563274374
$ws.Range("V8").Value = 41.637031421196575

# Row 9 (sub_leaf_8) - 6 cross angles (R,S,T,U,V,W)
$ws.Range("C9").Value = -13.420640582687634
$ws.Range("D9").Value = -12.370673454746896
$ws.Range("E9").Value = 48.18750805764982
$ws.Range("R9").Value = 30.72159022740245
$ws.Range("S9").Value = -71.62354316451218
$ws.Range("T9").Value = -52.75232665855775
$ws.Range("U9").Value = 28.010979749063956
$ws.Range("V9").Value = -58.72872507278277
$ws.Range("W9").Value = 43.848181423260485

# Row 10 (sub_leaf_9) - 5 cross angles (R,S,T,U,V)
$ws.Range("C10").Value = -24.298569096677017
$ws.Range("D10").Value = -47.560640312697444
$ws.Range("E10").Value = 51.51886498524703
$ws.Range("R10").Value = -92.16107565953361
$ws.Range("S10").Value = -53.60667365291955
$ws.Range("T10").Value = 38.5199838889373
$ws.Range("U10").Value = -47.560640312697444
$ws.Range("V10").Value = 33.31556025282822

# Row 11 (sub_leaf_10) - 4 cross angles (R,S,T,U)
$ws.Range("C11").Value = 0.5603904517009131
$ws.Range("D11").Value = 1.2049145288685565
$ws.Range("E11").Value = 39.43792487562386
$ws.Range("R11").Value = -41.67196532688721
$ws.Range("S11").Value = 38.35758965183823
$ws.Range("T11").Value = 41.50369807595375
$ws.Range("U11").Value = -35.94776059410112
